$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Nathan Grey row (row 3 of the table; row 1 is the header) - R3/R4 column (col 4)
$cellNathan = $t.Cell(3, 4)
$cellNathan.Range.InsertAfter("Loadr3()`ryield()`rdocumentation")

# Bradley Kersting row (row 4) - R3/R4 column (col 4)
$cellBradley = $t.Cell(4, 4)
$cellBradley.Range.InsertAfter("sys_call`ralarm process`rdocumentation")

# Lennon Jones row (row 5, last row of the table) - R3/R4 column (col 4)
$cellLennon = $t.Cell(5, 4)
$cellLennon.Range.InsertAfter("Alarm process`rTesting`rcomhand process`rMeeting with Sam")
